$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (LEIRIA / PENICHE) with the refreshed election results
$ws.Range("H2").Value  = 268
$ws.Range("I2").Value  = 717
$ws.Range("J2").Value  = 3006
$ws.Range("K2").Value  = 13
$ws.Range("L2").Value  = 839
$ws.Range("M2").Value  = 42
$ws.Range("N2").Value  = 494
$ws.Range("P2").Value  = 12
$ws.Range("Q2").Value  = 2
$ws.Range("R2").Value  = 44
$ws.Range("S2").Value  = 333
$ws.Range("T2").Value  = 491
$ws.Range("U2").Value  = 40
$ws.Range("V2").Value  = 4627
$ws.Range("W2").Value  = 2
$ws.Range("X2").Value  = 4740
$ws.Range("Y2").Value  = 6
$ws.Range("Z2").Value  = 71
$ws.Range("AA2").Value = 26
